$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Split Map")

# Consolidated colors / reordered inputs on acoustic setup (row 3 relabeling)
$ws.Range("A3").Value = "DSR Vocals"
$ws.Range("B3").Value = "DSR Bass"
$ws.Range("C3").Value = "DSR Aux"
$ws.Range("D3").Value = "DCR Vocals"
$ws.Range("E3").Value = "DCR Acoustic"
$ws.Range("F3").Value = "DCL Vocals"
$ws.Range("G3").Value = "DCL Electric"
$ws.Range("H3").Value = "DSL Vocals"
$ws.Range("J3").Value = "DSL Piano"
$ws.Range("K3").Value = "DSL Keys"
$ws.Range("L3").Value = "DSL Electric"
$ws.Range("M3").Value = "DSL Acoustic"
$ws.Range("P3").Value = "UC Drum Overhd"
$ws.Range("Q3").Value = "UC Drum Overhd"

# Row 7 relabel
$ws.Range("A7").Value = "UC Vocals"

# Fixed position / selection
$ws.Range("A8").Select()
